# Updates cryptos list figures (price + 1h volume change) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
# Numeric-looking price strings (e.g. "1.00", "0.0789") are written via a
# Formula->Copy->PasteSpecial(values) round-trip so they land back in the
# sheet as literal text (matching the original inlineStr cells) instead of
# being auto-coerced to a Number by plain .Value assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.001.16'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '2.260.85'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Formula = '="305.28"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Formula = '="95.42"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +2.94%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Formula = '="0.489"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +0.83%  '
$ws.Range('D10').Formula = '="35.00"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +7.37%  '
$ws.Range('D11').Formula = '="0.0789"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Formula = '="6.62"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = '2.615.59'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Formula = '="14.35"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '2.229.22'
$ws.Range('E16').Value = '  -1.75%  '
$ws.Range('D17').Formula = '="0.791"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '41.921.88'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').Formula = '="12.36"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -4.40%  '
$ws.Range('D20').Value = '0.0₃0901'
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('D21').Formula = '="5.96"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').Formula = '="67.56"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('D23').Formula = '="237.45"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('D24').Formula = '="2.57"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Formula = '="1.00"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Formula = '="1.92"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').Formula = '="23.69"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('D28').Formula = '="36.56"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  +5.49%  '
$ws.Range('D29').Formula = '="9.49"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('D30').Formula = '="2.10"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +1.63%  '
$ws.Range('D31').Formula = '="160.12"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('D32').Formula = '="5.21"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  -2.45%  '
$ws.Range('D33').Formula = '="1.00"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Formula = '="3.17"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('D35').Formula = '="0.0737"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').Formula = '="17.01"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').Formula = '="1.82"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('D41').Formula = '="3.99"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('E42').Value = '  +6.78%  '
$ws.Range('D43').Value = '1.980.39'
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').Formula = '="18.91"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -3.93%  '
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').Formula = '="2.92"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('D47').Formula = '="9.91"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -3.40%  '
$ws.Range('D48').Formula = '="53.13"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('D49').Formula = '="72.25"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Formula = '="90.84"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  -0.87%  '

$excel.CutCopyMode = $false
